$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "The live axel design..." paragraph: wrap "drives" in a gramStart/gramEnd
#    proofErr pair (splitting the single run into three runs).
# ---------------------------------------------------------------------------
$findRange1 = $d.Content
$ok1 = $findRange1.Find.Execute(
    "The live axel design is driven by a main motor and steered my changing the angle of wheels, much like how a car drives. This uses the least motors, reducing cost, but is the only option that cannot turn in place.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($ok1) {
    $target1 = $d.Range($findRange1.Start, $findRange1.End)
    $xml1 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The live axel design is driven by a main motor and steered my changing the angle of wheels, much like how a car </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>drives</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>. This uses the least motors, reducing cost, but is the only option that cannot turn in place.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target1.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# 2) Remove the stray "_GoBack" bookmark from the end of the document (after
#    "...notoriously durable, even being used in body armor.") — it moves to
#    the location of the most recent edit (see step 3 below).
# ---------------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
}

# ---------------------------------------------------------------------------
# 3) "Rollers are use a studded cylinder..." -> "Rollers use a studded
#    cylinder...", and drop a fresh "_GoBack" bookmark right after "Rollers"
#    (marking this as the most recent edit location).
# ---------------------------------------------------------------------------
$findRange2 = $d.Content
$ok2 = $findRange2.Find.Execute("Rollers are use", $true, $false, $false, $false, $false, $true, 1, $false, "Rollers use", 2)

$findRange3 = $d.Content
$ok3 = $findRange3.Find.Execute("Rollers", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok3) {
    $afterRollers = $d.Range($findRange3.End, $findRange3.End)
    $d.Bookmarks.Add("_GoBack", $afterRollers)
}
